# Insert a new price-report row for Papaya ("Vega Modelo de Temuco")
# immediately above the existing row 26, pushing every subsequent row
# (old 26..82) down by one (new 27..83). The worksheet's used-range
# dimension grows from A1:T82 to A1:T83 automatically on insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 26:82 down to 27:83, leaving a blank row 26 in their place.
$ws.Rows.Item(26).Insert()

# Populate the newly-inserted row 26 with the new weekly record.
$ws.Range("A26").Value = 10
$ws.Range("B26").Value = "Vega Modelo de Temuco"
$ws.Range("C26").Value = "La Araucanía"
$ws.Range("D26").Value = 44965
$ws.Range("E26").Value = 9
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100108
$ws.Range("H26").Value = "Tropicales y subtropicales"
$ws.Range("I26").Value = 100108004
$ws.Range("J26").Value = "Papaya"
$ws.Range("K26").Value = "Cultivar IV Región"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 40
$ws.Range("N26").Value = 32000
$ws.Range("O26").Value = 32000
$ws.Range("P26").Value = 32000
$ws.Range("Q26").Value = "`$/bandeja 10 kilos"
$ws.Range("R26").Value = "Provincia del Elquí"
$ws.Range("S26").Value = 3200
$ws.Range("T26").Value = 10
